# Updated cryptos list on Thu May 23 21:34:24 UTC 2024 with GitHub Actions
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) values for
# each coin row. Values that look like plain numbers (e.g. "0.999",
# "594.62") are written with a leading apostrophe so Excel keeps them as
# text (matching the original inlineStr cells) instead of converting them
# to numeric values; the cell style is then reset to "Normal" so no stray
# quote-prefix formatting is left on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "67.428.64"
$ws.Cells.Item(2, 5).Value = "  -2.75%  "

$ws.Cells.Item(3, 4).Value = "3.779.99"
$ws.Cells.Item(3, 5).Value = "  +1.06%  "

$ws.Cells.Item(4, 4).Value = "'0.999"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = "  -0.06%  "

$ws.Cells.Item(5, 4).Value = "'594.62"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -3.27%  "

$ws.Cells.Item(6, 4).Value = "'175.72"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -1.49%  "

$ws.Cells.Item(7, 4).Value = "3.779.39"
$ws.Cells.Item(7, 5).Value = "  +1.09%  "

$ws.Cells.Item(8, 4).Value = "'0.999"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  -0.18%  "

$ws.Cells.Item(9, 4).Value = "'0.519"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  -1.29%  "

$ws.Cells.Item(10, 4).Value = "'0.159"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  -4.33%  "

$ws.Cells.Item(11, 5).Value = "  -5.77%  "

$ws.Cells.Item(12, 5).Value = "  -3.94%  "

$ws.Cells.Item(13, 4).Value = "'38.08"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  -4.59%  "

$ws.Cells.Item(14, 4).Value = "'0.0000243"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  -3.82%  "

$ws.Cells.Item(15, 4).Value = "4.412.93"
$ws.Cells.Item(15, 5).Value = "  +1.15%  "

$ws.Cells.Item(16, 4).Value = "3.784.20"
$ws.Cells.Item(16, 5).Value = "  +1.21%  "

$ws.Cells.Item(17, 4).Value = "67.451.83"
$ws.Cells.Item(17, 5).Value = "  -2.82%  "

$ws.Cells.Item(18, 5).Value = "  -4.67%  "

$ws.Cells.Item(19, 5).Value = "  -3.65%  "

$ws.Cells.Item(20, 4).Value = "'16.32"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  -0.10%  "

$ws.Cells.Item(21, 4).Value = "'487.94"
$ws.Cells.Item(21, 4).Style = "Normal"

$ws.Cells.Item(22, 4).Value = "'9.04"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  -1.46%  "

$ws.Cells.Item(23, 4).Value = "'0.724"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +0.45%  "

$ws.Cells.Item(24, 4).Value = "'0.0000148"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +10.09%  "

$ws.Cells.Item(25, 4).Value = "'83.88"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -2.00%  "

$ws.Cells.Item(26, 4).Value = "'2.34"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -9.59%  "

$ws.Cells.Item(27, 5).Value = "  -5.52%  "

$ws.Cells.Item(28, 4).Value = "'10.18"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  -7.52%  "

$ws.Cells.Item(29, 5).Value = "  +0.15%  "

$ws.Cells.Item(30, 4).Value = "'2.93"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +0.57%  "

$ws.Cells.Item(31, 5).Value = "  -2.69%  "

$ws.Cells.Item(32, 4).Value = "'32.82"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +8.18%  "

$ws.Cells.Item(33, 4).Value = "'7.73"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -4.04%  "

$ws.Cells.Item(34, 4).Value = "'0.108"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  -3.81%  "

$ws.Cells.Item(35, 4).Value = "'0.999"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +0.05%  "

$ws.Cells.Item(36, 5).Value = "  -3.60%  "

$ws.Cells.Item(37, 4).Value = "'0.135"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -0.93%  "

$ws.Cells.Item(38, 4).Value = "'5.75"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  -5.91%  "

$ws.Cells.Item(39, 4).Value = "'0.327"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  -6.22%  "

$ws.Cells.Item(40, 4).Value = "'452.97"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +0.68%  "

$ws.Cells.Item(41, 4).Value = "'49.02"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  -1.29%  "

$ws.Cells.Item(42, 4).Value = "'1.99"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -3.39%  "

$ws.Cells.Item(43, 5).Value = "  -7.37%  "

$ws.Cells.Item(44, 4).Value = "'8.31"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  -2.83%  "

$ws.Cells.Item(45, 4).Value = "'41.22"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -7.56%  "

$ws.Cells.Item(46, 4).Value = "2.814.41"
$ws.Cells.Item(46, 5).Value = "  -4.44%  "

$ws.Cells.Item(47, 4).Value = "'141.46"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +2.07%  "

$ws.Cells.Item(48, 5).Value = "  +0.02%  "

$ws.Cells.Item(49, 4).Value = "'0.0348"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  -3.10%  "

$ws.Cells.Item(50, 4).Value = "'25.86"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -4.76%  "

$ws.Cells.Item(51, 4).Value = "'23.20"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +7.38%  "
